# Update crypto price/volume data in the active sheet (Price = column D, Volume(1h) = column E).
# Source values are text (e.g. "299.51", "-1.94%"), so force the cells to remain
# text cells (NumberFormat "@") before writing, otherwise Excel auto-coerces
# numeric-looking strings (and "%" strings) into actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "299.51";        E = "-1.94%" },
    @{ Row = 3;  D = "31.32";         E = "-1.34%" },
    @{ Row = 4;  D = "5.106";         E = "-1.59%" },
    @{ Row = 5;  D = "0.07964";       E = "5.59%" },
    @{ Row = 6;  D = "2.352";         E = "0.95%" },
    @{ Row = 7;  D = "7.769";         E = "-3.04%" },
    @{ Row = 8;  D = "3.863";         E = "-0.45%" },
    @{ Row = 9;  D = "0.9228";        E = "0.85%" },
    @{ Row = 10; D = "0.1736";        E = "-0.60%" },
    @{ Row = 11; D = "0.07583";       E = "0.40%" },
    @{ Row = 12; D = "0.09342";       E = "13.06%" },
    @{ Row = 13; E = "0.45%" },
    @{ Row = 14; D = "0.1003";        E = "0.84%" },
    @{ Row = 15; D = "0.001512";      E = "0.35%" },
    @{ Row = 16; D = "0.006022";      E = "-1.44%" },
    @{ Row = 17; D = "3.485";         E = "-0.49%" },
    @{ Row = 18; D = "2.268";         E = "1.35%" },
    @{ Row = 20; E = "-0.24%" },
    @{ Row = 21; D = "3.925";         E = "-15.66%" },
    @{ Row = 23; D = "0.04631";       E = "0.14%" },
    @{ Row = 24; D = "0.001250" },
    @{ Row = 25; D = "0.004481";      E = "-1.27%" },
    @{ Row = 26; D = "0.0001200";     E = "-7.51%" },
    @{ Row = 27; D = "0.0003396" },
    @{ Row = 39; D = "0.01756";       E = "-1.69%" },
    @{ Row = 40; D = "0.04626";       E = "0.58%" },
    @{ Row = 41; D = "0.006975";      E = "-4.50%" },
    @{ Row = 42; E = "-0.31%" },
    @{ Row = 43; D = "0.002189";      E = "-2.04%" },
    @{ Row = 44; D = "0.01028";       E = "-4.75%" },
    @{ Row = 45; D = "0.00006276";    E = "-3.36%" },
    @{ Row = 46; D = "0.00000000750"; E = "-0.04%" },
    @{ Row = 47; D = "0.007977";      E = "-19.24%" },
    @{ Row = 48; D = "0.7467";        E = "-9.00%" },
    @{ Row = 49; D = "0.00002099";    E = "-0.04%" },
    @{ Row = 50; D = "0.0001999";     E = "0.03%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}
